# Daily "GitHub Actions" refresh of the cryptos price table.
# Price/Volume(1h) text for most rows is updated in place; three coin pairs
# (rows 18/19, 29/30, 35/36) swap rank position, so those rows get their
# Coin/Link/Price/Volume cells fully rewritten rather than just Price/Volume.
#
# Cells are written via .Value2 (not .Value) so Excel's COM layer never
# coerces a text cell to a different OLE variant type. Column D holds prices
# formatted as plain text (e.g. "8.30", "1.00", "51.131.72" with multiple
# dots as thousands separators) - for any value that is *also* syntactically
# a valid number, a leading apostrophe is prefixed (mirrors a user typing
# '8.30 into the cell) so Excel stores the literal digits/trailing zeros
# instead of silently re-parsing it as a float and dropping them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '50.910.92'
$ws.Range("E2").Value2 = '  -0.45%  '

$ws.Range("D3").Value2 = '2.937.36'
$ws.Range("E3").Value2 = '  -0.65%  '

$ws.Range("E4").Value2 = '  +0.23%  '

$ws.Range("D5").Value2 = '''374.76'
$ws.Range("E5").Value2 = '  -1.32%  '

$ws.Range("D6").Value2 = '''101.72'
$ws.Range("E6").Value2 = '  -2.94%  '

$ws.Range("D7").Value2 = '''0.535'
$ws.Range("E7").Value2 = '  -0.88%  '

$ws.Range("D9").Value2 = '''0.584'
$ws.Range("E9").Value2 = '  -1.82%  '

$ws.Range("D10").Value2 = '''36.27'
$ws.Range("E10").Value2 = '  -2.14%  '

$ws.Range("E11").Value2 = '  -0.52%  '

$ws.Range("D12").Value2 = '''0.0837'
$ws.Range("E12").Value2 = '  -0.29%  '

$ws.Range("D13").Value2 = '3.415.24'
$ws.Range("E13").Value2 = '  -0.07%  '

$ws.Range("D14").Value2 = '''17.88'
$ws.Range("E14").Value2 = '  -2.97%  '

$ws.Range("D15").Value2 = '''7.39'
$ws.Range("E15").Value2 = '  -1.37%  '

$ws.Range("D16").Value2 = '2.937.55'
$ws.Range("E16").Value2 = '  -0.56%  '

$ws.Range("D17").Value2 = '''0.981'
$ws.Range("E17").Value2 = '  +1.71%  '

$ws.Range("B18").Value2 = 'Uniswap'
$ws.Range("C18").Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value2 = '''10.52'
$ws.Range("E18").Value2 = '  +42.68%  '

$ws.Range("B19").Value2 = 'WrappedBTC'
$ws.Range("C19").Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value2 = '50.953.76'
$ws.Range("E19").Value2 = '  -0.41%  '

$ws.Range("D20").Value2 = '''3.14'
$ws.Range("E20").Value2 = '  -5.55%  '

$ws.Range("D21").Value2 = '''12.61'
$ws.Range("E21").Value2 = '  -1.88%  '

$ws.Range("D22").Value2 = '0.0₃0954'
$ws.Range("E22").Value2 = '  -0.63%  '

$ws.Range("D23").Value2 = '''263.92'
$ws.Range("E23").Value2 = '  +1.01%  '

$ws.Range("D24").Value2 = '''68.35'
$ws.Range("E24").Value2 = '  -1.59%  '

$ws.Range("D25").Value2 = '''2.99'
$ws.Range("E25").Value2 = '  +5.98%  '

$ws.Range("D26").Value2 = '''8.30'
$ws.Range("E26").Value2 = '  +8.42%  '

$ws.Range("D27").Value2 = '''7.92'
$ws.Range("E27").Value2 = '  +7.84%  '

$ws.Range("D28").Value2 = '''0.169'
$ws.Range("E28").Value2 = '  -0.25%  '

$ws.Range("B29").Value2 = 'Dai'
$ws.Range("C29").Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").Value2 = '''1.00'
$ws.Range("E29").Value2 = '  -0.10%  '

$ws.Range("B30").Value2 = 'Hedera'
$ws.Range("C30").Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value2 = '''0.112'
$ws.Range("E30").Value2 = '  +0.80%  '

$ws.Range("D31").Value2 = '''25.57'
$ws.Range("E31").Value2 = '  -0.98%  '

$ws.Range("D32").Value2 = '''9.86'
$ws.Range("E32").Value2 = '  +0.26%  '

$ws.Range("D33").Value2 = '''50.67'
$ws.Range("E33").Value2 = '  -0.58%  '

$ws.Range("D34").Value2 = '''33.56'
$ws.Range("E34").Value2 = '  -3.29%  '

$ws.Range("B35").Value2 = 'Toncoin'
$ws.Range("C35").Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D35").Value2 = '''2.02'
$ws.Range("E35").Value2 = '  -2.98%  '

$ws.Range("B36").Value2 = 'VeChain'
$ws.Range("C36").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value2 = '''0.0443'
$ws.Range("E36").Value2 = '  -1.03%  '

$ws.Range("E37").Value2 = '  +0.11%  '

$ws.Range("D38").Value2 = '''2.99'
$ws.Range("E38").Value2 = '  -2.34%  '

$ws.Range("D39").Value2 = '''2.54'
$ws.Range("E39").Value2 = '  -1.19%  '

$ws.Range("D40").Value2 = '''0.115'
$ws.Range("E40").Value2 = '  -0.35%  '

$ws.Range("D41").Value2 = '''16.37'
$ws.Range("E41").Value2 = '  -5.13%  '

$ws.Range("E42").Value2 = '  -3.03%  '

$ws.Range("D43").Value2 = '''120.36'
$ws.Range("E43").Value2 = '  -2.51%  '

$ws.Range("D44").Value2 = '''0.293'
$ws.Range("E44").Value2 = '  +0.69%  '

$ws.Range("D45").Value2 = '''20.94'
$ws.Range("E45").Value2 = '  -4.79%  '

$ws.Range("D46").Value2 = '''2.02'
$ws.Range("E46").Value2 = '  -1.97%  '

$ws.Range("D47").Value2 = '''3.29'
$ws.Range("E47").Value2 = '  +2.26%  '

$ws.Range("E48").Value2 = '  -3.11%  '

$ws.Range("D49").Value2 = '1.986.76'

$ws.Range("D50").Value2 = '''0.0342'
$ws.Range("E50").Value2 = '  -1.42%  '

$ws.Range("B51").Value2 = 'TrustWalletToken'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value2 = '''1.27'
$ws.Range("E51").Value2 = '  -1.64%  '

